# Apply updated figures to the analytics workbook.
# This updates individual data cells across the four percentage/count
# sheets ("10 Year Percentages", "10 Year Counts", "5 Year Percentages",
# "5 Year Counts") to reflect the latest computed values.

$wb = $excel.ActiveWorkbook

# --- 10 Year Percentages ---
$ws = $wb.Worksheets.Item("10 Year Percentages")
$ws.Range("B5").Value = 60
$ws.Range("D6").Value = 65
$ws.Range("E6").Value = 71
$ws.Range("B7").Value = 56.00000000000001
$ws.Range("C7").Value = 61
$ws.Range("D7").Value = 70
$ws.Range("E7").Value = 77
$ws.Range("B8").Value = 66
$ws.Range("C8").Value = 64
$ws.Range("F8").Value = 71
$ws.Range("B9").Value = 46
$ws.Range("E9").Value = 76
$ws.Range("F9").Value = 74
$ws.Range("D10").Value = 68
$ws.Range("D12").Value = 68
$ws.Range("C13").Value = 71

# --- 10 Year Counts ---
$ws = $wb.Worksheets.Item("10 Year Counts")
$ws.Range("B5").Value = 15
$ws.Range("D6").Value = 94
$ws.Range("E6").Value = 28
$ws.Range("B7").Value = 36
$ws.Range("C7").Value = 148
$ws.Range("D7").Value = 88
$ws.Range("E7").Value = 26
$ws.Range("B8").Value = 79
$ws.Range("C8").Value = 332
$ws.Range("D8").Value = 214
$ws.Range("F8").Value = 35
$ws.Range("B9").Value = 35
$ws.Range("C9").Value = 188
$ws.Range("D9").Value = 104
$ws.Range("E9").Value = 38
$ws.Range("F9").Value = 19
$ws.Range("C10").Value = 254
$ws.Range("D10").Value = 171
$ws.Range("D12").Value = 63
$ws.Range("C13").Value = 56

# --- 5 Year Percentages ---
$ws = $wb.Worksheets.Item("5 Year Percentages")
$ws.Range("B3").Value = 50
$ws.Range("D4").Value = 59
$ws.Range("E4").Value = 70
$ws.Range("B5").Value = 54
$ws.Range("C5").Value = 61
$ws.Range("D5").Value = 70
$ws.Range("B6").Value = 57.99999999999999
$ws.Range("C6").Value = 64
$ws.Range("D6").Value = 64
$ws.Range("F6").Value = 75
$ws.Range("B7").Value = 56.99999999999999
$ws.Range("C7").Value = 57.99999999999999
$ws.Range("D7").Value = 64
$ws.Range("E7").Value = 76
$ws.Range("F7").Value = 60
$ws.Range("D8").Value = 63
$ws.Range("D10").Value = 84
$ws.Range("C11").Value = 61

# --- 5 Year Counts ---
$ws = $wb.Worksheets.Item("5 Year Counts")
$ws.Range("B3").Value = 4
$ws.Range("D4").Value = 29
$ws.Range("E4").Value = 10
$ws.Range("B5").Value = 13
$ws.Range("C5").Value = 41
$ws.Range("D5").Value = 30
$ws.Range("E5").Value = 9
$ws.Range("B6").Value = 26
$ws.Range("C6").Value = 102
$ws.Range("D6").Value = 78
$ws.Range("F6").Value = 16
$ws.Range("B7").Value = 14
$ws.Range("C7").Value = 60
$ws.Range("D7").Value = 44
$ws.Range("E7").Value = 21
$ws.Range("F7").Value = 10
$ws.Range("C8").Value = 95
$ws.Range("D8").Value = 71
$ws.Range("D10").Value = 31
$ws.Range("C11").Value = 28
